$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header style (from H1) onto I1:J1,
#     then set the header labels.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (2-16): plain numeric values in columns I and J.
$values = @{
    2  = @(6, 6)
    3  = @(8, 8)
    4  = @(8, 8)
    5  = @(6, 6)
    6  = @(6, 7)
    7  = @(8, 8)
    8  = @(10, 10)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(6, 6)
    12 = @(9, 9)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
